# Applies the "SLL64 errors fixed (not tested)" activity-log entries to the
# "Activity Log - Part 3" worksheet: fills in the student header info and the
# first several log rows (date/start/end time + description) describing work
# on reviewing Part 2 and implementing/debugging the Barrel Shifter / MUX.

$wb = $excel.ActiveWorkbook

# Tidy up the view on the other two log sheets (scroll back to top / move the
# selection) before doing the real data entry on Part 3.
$wsPart1 = $wb.Worksheets.Item("Activity Log - Part 1")
$wsPart1.Activate()
$wsPart1.Range("G28").Select()

$wsPart2 = $wb.Worksheets.Item("Activity Log - Part 2")
$wsPart2.Activate()
$wsPart2.Range("A6:G11").Select()

$ws = $wb.Worksheets.Item("Activity Log - Part 3")
$ws.Activate()

# --- Header block -----------------------------------------------------
$ws.Range("B1").Value = "Ruelt Yean (Ryan), Kiew"
$ws.Range("B2").Value = 301290779
$ws.Range("B3").Value = "G47"

# --- Activity rows ------------------------------------------------------
# Row 6
$ws.Range("B6").Value = 779
$ws.Range("C6").Value = 43931
$ws.Range("D6").Value = 0.64513888888888882
$ws.Range("E6").Value = 0.65555555555555556
$ws.Range("G6").Value = "Read through Part 2 pdf"

# Row 7
$ws.Range("B7").Value = 779
$ws.Range("C7").Value = 43931
$ws.Range("D7").Value = 0.65555555555555556
$ws.Range("E7").Value = 0.66180555555555554
$ws.Range("G7").Value = "Reviewed Barrel Shifter design implementation"

# Row 8
$ws.Range("C8").Value = 43931
$ws.Range("D8").Value = 0.72638888888888886
$ws.Range("E8").Value = 0.75138888888888899
$ws.Range("G8").Value = "First implementation  of Barrel Shifter"

# Row 9
$ws.Range("C9").Value = 43931
$ws.Range("D9").Value = 0.84097222222222223
$ws.Range("E9").Value = 0.94027777777777777
$ws.Range("G9").Value = "Second implementation of Barrel Shifter; attempting to use same MUX entity"

# Row 10
$ws.Range("C10").Value = 43931
$ws.Range("D10").Value = 0.58958333333333335
$ws.Range("E10").Value = 0.6069444444444444

# Row 11
$ws.Range("D11").Value = 0.83611111111111114
$ws.Range("E11").Value = 0.96111111111111114

# Row 12
$ws.Range("D12").Value = 0.4381944444444445
$ws.Range("E12").Value = 0.47291666666666665
$ws.Range("G12").Value = "Third implementaiton of Barrel Shifter; error fixing, split MUX into three entities "
